$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$row = 52

# Column A holds a date-like string ("2018.08.27"). Excel's automatic data
# detection would otherwise convert it into a real date serial number, so we
# temporarily force a text number format while the value is entered, then
# restore the cell's style back to the default "Normal" style so the cell
# keeps no explicit style index (matching the rest of the sheet).
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2018.08.27"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = "16:25:47"
$ws.Cells.Item($row, 3).Value = "RS"
$ws.Cells.Item($row, 4).Value = 10
$ws.Cells.Item($row, 5).Value = 250
$ws.Cells.Item($row, 6).Value = 0.1
$ws.Cells.Item($row, 7).Value = 0.96
$ws.Cells.Item($row, 8).Value = 2975
$ws.Cells.Item($row, 9).Value = 0.43
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = "N/A"
$ws.Cells.Item($row, 12).Value = "N/A"
